# Apply re-sorted match data for Slovakia 2-liga 2023-2024 (script run 05-11-2023 20:45)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 105
$ws.Cells.Item(105, 6).Value = "Puchov"
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = "Pohronie"
$ws.Cells.Item(105, 9).Value = 4
$ws.Cells.Item(105, 10).Value = 1.75
$ws.Cells.Item(105, 12).Value = 1.71
$ws.Cells.Item(105, 13).Value = "28/10/2023 14:21"
$ws.Cells.Item(105, 14).Value = 3.6
$ws.Cells.Item(105, 16).Value = 3.81
$ws.Cells.Item(105, 17).Value = "28/10/2023 14:21"
$ws.Cells.Item(105, 18).Value = 3.77
$ws.Cells.Item(105, 20).Value = 4.45
$ws.Cells.Item(105, 21).Value = "28/10/2023 14:21"
$ws.Cells.Item(105, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/msk-puchov-fk-pohronie/QDGXEwbG/"

# Row 106
$ws.Cells.Item(106, 6).Value = "Trebisov"
$ws.Cells.Item(106, 7).Value = 2
$ws.Cells.Item(106, 8).Value = "Zilina B"
$ws.Cells.Item(106, 9).Value = 2
$ws.Cells.Item(106, 10).Value = 2.77
$ws.Cells.Item(106, 12).Value = 2.31
$ws.Cells.Item(106, 13).Value = "28/10/2023 14:15"
$ws.Cells.Item(106, 14).Value = 3.43
$ws.Cells.Item(106, 16).Value = 3.69
$ws.Cells.Item(106, 17).Value = "28/10/2023 14:16"
$ws.Cells.Item(106, 18).Value = 2.16
$ws.Cells.Item(106, 20).Value = 2.68
$ws.Cells.Item(106, 21).Value = "28/10/2023 14:15"
$ws.Cells.Item(106, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/trebisov-zilina/rkRSFJrA/"

# Row 108
$ws.Cells.Item(108, 6).Value = "Povazska Bystrica"
$ws.Cells.Item(108, 7).Value = 3
$ws.Cells.Item(108, 8).Value = "Slovan Bratislava B"
$ws.Cells.Item(108, 9).Value = 1
$ws.Cells.Item(108, 10).Value = 1.8
$ws.Cells.Item(108, 12).Value = 1.46
$ws.Cells.Item(108, 13).Value = "28/10/2023 13:57"
$ws.Cells.Item(108, 14).Value = 3.56
$ws.Cells.Item(108, 16).Value = 4.44
$ws.Cells.Item(108, 17).Value = "28/10/2023 13:57"
$ws.Cells.Item(108, 18).Value = 3.58
$ws.Cells.Item(108, 20).Value = 6.15
$ws.Cells.Item(108, 21).Value = "28/10/2023 13:57"
$ws.Cells.Item(108, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/povazska-bystrica-slovan-bratislava/ObZzEcDM/"

# Row 109
$ws.Cells.Item(109, 6).Value = "Samorin"
$ws.Cells.Item(109, 7).Value = 2
$ws.Cells.Item(109, 8).Value = "Myjava"
$ws.Cells.Item(109, 9).Value = 2
$ws.Cells.Item(109, 10).Value = 2.32
$ws.Cells.Item(109, 12).Value = 2.78
$ws.Cells.Item(109, 13).Value = "29/10/2023 09:48"
$ws.Cells.Item(109, 14).Value = 3.35
$ws.Cells.Item(109, 16).Value = 3.53
$ws.Cells.Item(109, 17).Value = "29/10/2023 09:48"
$ws.Cells.Item(109, 18).Value = 2.66
$ws.Cells.Item(109, 20).Value = 2.33
$ws.Cells.Item(109, 21).Value = "29/10/2023 09:48"
$ws.Cells.Item(109, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/samorin-myjava/Ua2cVbrc/"

# Row 111
$ws.Cells.Item(111, 6).Value = "Petrzalka"
$ws.Cells.Item(111, 7).Value = 3
$ws.Cells.Item(111, 8).Value = "Komarno"
$ws.Cells.Item(111, 9).Value = 1
$ws.Cells.Item(111, 10).Value = 2.89
$ws.Cells.Item(111, 12).Value = 2.91
$ws.Cells.Item(111, 13).Value = "29/10/2023 10:02"
$ws.Cells.Item(111, 14).Value = 3.2
$ws.Cells.Item(111, 16).Value = 3.44
$ws.Cells.Item(111, 17).Value = "29/10/2023 10:02"
$ws.Cells.Item(111, 18).Value = 2.2
$ws.Cells.Item(111, 20).Value = 2.29
$ws.Cells.Item(111, 21).Value = "29/10/2023 10:02"
$ws.Cells.Item(111, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/petrzalka-komarno/GdEF6F5d/"

# Row 117
$ws.Cells.Item(117, 6).Value = "Slovan Bratislava B"
$ws.Cells.Item(117, 7).Value = 2
$ws.Cells.Item(117, 8).Value = "Puchov"
$ws.Cells.Item(117, 9).Value = 4
$ws.Cells.Item(117, 10).Value = 2.59
$ws.Cells.Item(117, 11).Value = "03/11/2023 22:44"
$ws.Cells.Item(117, 12).Value = 2.59
$ws.Cells.Item(117, 13).Value = "03/11/2023 22:44"
$ws.Cells.Item(117, 14).Value = 3.43
$ws.Cells.Item(117, 15).Value = "05/11/2023 08:33"
$ws.Cells.Item(117, 16).Value = 3.43
$ws.Cells.Item(117, 17).Value = "05/11/2023 08:33"
$ws.Cells.Item(117, 18).Value = 2.4
$ws.Cells.Item(117, 19).Value = "03/11/2023 22:44"
$ws.Cells.Item(117, 20).Value = 2.4
$ws.Cells.Item(117, 21).Value = "03/11/2023 22:44"
$ws.Cells.Item(117, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/slovan-bratislava-msk-puchov/AJZlLEKk/"

# Row 118
$ws.Cells.Item(118, 6).Value = "Malzenice"
$ws.Cells.Item(118, 7).Value = 4
$ws.Cells.Item(118, 8).Value = "Samorin"
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 10).Value = 2.06
$ws.Cells.Item(118, 11).Value = "05/11/2023 08:04"
$ws.Cells.Item(118, 12).Value = 2.06
$ws.Cells.Item(118, 13).Value = "05/11/2023 08:04"
$ws.Cells.Item(118, 14).Value = 3.52
$ws.Cells.Item(118, 15).Value = "05/11/2023 08:30"
$ws.Cells.Item(118, 16).Value = 3.52
$ws.Cells.Item(118, 17).Value = "05/11/2023 08:30"
$ws.Cells.Item(118, 18).Value = 3.27
$ws.Cells.Item(118, 19).Value = "05/11/2023 08:04"
$ws.Cells.Item(118, 20).Value = 3.27
$ws.Cells.Item(118, 21).Value = "05/11/2023 08:04"
$ws.Cells.Item(118, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/malzenice-samorin/4tZhKYZe/"

# New row 119 (appended match record)
$ws.Cells.Item(118, 1).Copy($ws.Cells.Item(119, 1))
$ws.Cells.Item(118, 5).Copy($ws.Cells.Item(119, 5))

$ws.Cells.Item(119, 1).Value = 118
$ws.Cells.Item(119, 2).Value = "slovakia"
$ws.Cells.Item(119, 3).Value = "2-liga"
$ws.Cells.Item(119, 4).Value = "2023-2024"
$ws.Cells.Item(119, 5).Value = 45235.75
$ws.Cells.Item(119, 6).Value = "Myjava"
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = "Povazska Bystrica"
$ws.Cells.Item(119, 9).Value = 4
$ws.Cells.Item(119, 10).Value = 2.63
$ws.Cells.Item(119, 11).Value = "04/11/2023 06:12"
$ws.Cells.Item(119, 12).Value = 2.19
$ws.Cells.Item(119, 13).Value = "05/11/2023 17:57"
$ws.Cells.Item(119, 14).Value = 3.2
$ws.Cells.Item(119, 15).Value = "04/11/2023 06:12"
$ws.Cells.Item(119, 16).Value = 3.24
$ws.Cells.Item(119, 17).Value = "05/11/2023 17:57"
$ws.Cells.Item(119, 18).Value = 2.43
$ws.Cells.Item(119, 19).Value = "04/11/2023 06:12"
$ws.Cells.Item(119, 20).Value = 3.26
$ws.Cells.Item(119, 21).Value = "05/11/2023 17:58"
$ws.Cells.Item(119, 22).Value = "https://www.betexplorer.com/football/slovakia/2-liga/myjava-povazska-bystrica/QRe6TxD9/"
